$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 452, shifting existing rows 452-548 down to 453-549
$ws.Rows.Item(452).Insert()

# Populate the newly inserted row 452 with the new data record
$ws.Range("A452").Value = 10
$ws.Range("B452").Value = 'Vega Modelo de Temuco'
$ws.Range("C452").Value = 'La Araucanía'
$ws.Range("D452").Value = 45275
$ws.Range("E452").Value = 9
$ws.Range("F452").Value = 100112001
$ws.Range("G452").Value = 'Berenjena'
$ws.Range("H452").Value = 'Sin especificar'
$ws.Range("I452").Value = 'Primera'
$ws.Range("J452").Value = 35
$ws.Range("K452").Value = 14000
$ws.Range("L452").Value = 14000
$ws.Range("M452").Value = 14000
$ws.Range("N452").Value = '$/caja 40 unidades'
$ws.Range("O452").Value = 'Región de Arica y Parinacota'
$ws.Range("P452").Value = 350
$ws.Range("Q452").Value = 40
$ws.Range("R452").Value = 'Hortaliza'
